$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (ownTeam, oppTeam) before the current
# "batsman" column (D), shifting batsman..sr from D..I to F..K. ---
$ws.Range("D1:E1").EntireColumn.Insert()

# The numeric-looking columns (totalRuns, totalBalls, total4s, total6s, sr)
# are stored as text in this sheet (e.g. "187.50", "0.00"), so format those
# cells as Text *before* writing values, otherwise Excel would coerce them
# to real numbers and lose the formatting (leading/trailing zeros).
$ws.Range("G1:K4").NumberFormat = "@"

# --- Header row ---
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# --- Row 2 (existing match) - fill in new ownTeam/oppTeam values ---
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Royal Challengers Bangalore"

# --- Row 3 (new match: Abu Dhabi vs Mumbai Indians) ---
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 01 2020"
$ws.Range("C3").Value = "Mumbai won by 48 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Karun Nair "
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "3"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "0.00"

# --- Row 4 (new match: Dubai (DSC) vs Delhi Capitals) ---
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 20 2020"
$ws.Range("C4").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Delhi Capitals"
$ws.Range("F4").Value = "Karun Nair "
$ws.Range("G4").Value = "1"
$ws.Range("H4").Value = "3"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "33.33"
